$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "PASS" values that were in the L (STATUS) column for the three
# data rows -- the L column itself (header + width) stays intact, only the
# per-row "Header field" results are cleared out.
$ws.Range("L2:L4").ClearContents()

# Update the sheet's view: scroll so column G is left-most and select the
# (now empty) L2:L4 range as the active selection.
$win = $excel.ActiveWindow
$win.ScrollColumn = 7
$ws.Range("L2:L4").Select()
